$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (Coin name / Link / Volume%) - never numeric-looking, safe to assign directly.
function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell "D2" '43.643.92'
Set-TextCell "E2" '  +3.22%  '
Set-TextCell "D3" '2.318.74'
Set-TextCell "E3" '  +2.11%  '
Set-TextCell "E4" '  -0.02%  '
Set-TextCell "D5" '312.36'
Set-TextCell "E5" '  +1.75%  '
Set-TextCell "D6" '103.48'
Set-TextCell "E6" '  +6.59%  '
Set-TextCell "D7" '0.537'
Set-TextCell "E7" '  +2.09%  '
Set-TextCell "E8" '  -0.03%  '
Set-TextCell "D9" '0.536'
Set-TextCell "E9" '  +8.41%  '
Set-TextCell "D10" '36.12'
Set-TextCell "E10" '  +2.50%  '
Set-TextCell "E11" '  +3.65%  '
Set-TextCell "E12" '  -0.17%  '
Set-TextCell "E13" '  +2.27%  '
Set-TextCell "D14" '2.677.19'
Set-TextCell "E14" '  +2.08%  '
Set-TextCell "D15" '15.09'
Set-TextCell "E15" '  +2.37%  '
Set-TextCell "D16" '2.321.22'
Set-TextCell "E16" '  +2.18%  '
Set-TextCell "E17" '  +2.47%  '
Set-TextCell "D18" '43.542.51'
Set-TextCell "E18" '  +3.31%  '
Set-TextCell "D19" '12.55'
Set-TextCell "E19" '  +0.84%  '
Set-TextCell "D20" '0.0₃0934'
Set-TextCell "E20" '  +2.99%  '
Set-TextCell "E21" '  +2.27%  '
Set-TextCell "D22" '68.38'
Set-TextCell "E22" '  +0.35%  '
Set-TextCell "D23" '242.79'
Set-TextCell "E23" '  +1.81%  '
Set-TextCell "D24" '2.05'
Set-TextCell "E24" '  +5.70%  '
Set-TextCell "D25" '2.63'
Set-TextCell "E25" '  +2.28%  '
Set-TextCell "D26" '1.00'
Set-TextCell "E26" '  +0.06%  '
Set-TextCell "D27" '3.99'
Set-TextCell "E27" '  -1.49%  '
Set-TextCell "D28" '24.85'
Set-TextCell "E28" '  +5.14%  '
Set-TextCell "D29" '37.43'
Set-TextCell "E29" '  -1.04%  '
Set-TextCell "B30" 'Cosmos'
Set-TextCell "C30" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D30" '9.69'
Set-TextCell "E30" '  +1.87%  '
Set-TextCell "B31" 'Toncoin'
Set-TextCell "C31" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell "D31" '2.12'
Set-TextCell "E31" '  +0.10%  '
Set-TextCell "D32" '167.95'
Set-TextCell "E32" '  +4.10%  '
Set-TextCell "E33" '  +1.66%  '
Set-TextCell "E34" '  +0.01%  '
Set-TextCell "B35" 'LidoDAOToken'
Set-TextCell "C35" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell "D35" '3.13'
Set-TextCell "E35" '  -1.76%  '
Set-TextCell "B36" 'WEMIXToken'
Set-TextCell "C36" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell "D36" '2.53'
Set-TextCell "E36" '  +7.00%  '
Set-TextCell "D37" '0.0749'
Set-TextCell "E37" '  +1.58%  '
Set-TextCell "E38" '  +2.36%  '
Set-TextCell "E39" '  +2.01%  '
Set-TextCell "D40" '1.89'
Set-TextCell "E40" '  +3.13%  '
Set-TextCell "D41" '0.117'
Set-TextCell "E41" '  +2.05%  '
Set-TextCell "D42" '4.39'
Set-TextCell "E42" '  +8.21%  '
Set-TextCell "B43" 'ApeXProtocol'
Set-TextCell "C43" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D43" '2.32'
Set-TextCell "E43" '  -0.10%  '
Set-TextCell "B44" 'EnergySwap'
Set-TextCell "C44" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell "D44" '19.57'
Set-TextCell "E44" '  +3.88%  '
Set-TextCell "E45" '  +3.51%  '
Set-TextCell "D46" '1.977.56'
Set-TextCell "E46" '  +1.36%  '
Set-TextCell "E47" '  +4.59%  '
Set-TextCell "D48" '9.93'
Set-TextCell "D49" '55.84'
Set-TextCell "E49" '  +4.33%  '
Set-TextCell "D50" '2.94'
Set-TextCell "E50" '  +7.76%  '
Set-TextCell "D51" '1.58'
Set-TextCell "E51" '  +7.82%  '
